$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.800.73"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.270.97"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'304.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'92.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'32.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'53.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "'0.0796"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "2.623.10"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "'14.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "2.265.12"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "'0.778"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "41.721.28"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'12.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'5.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'67.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'243.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'23.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").Value = "'9.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'35.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("D32").Value = "'160.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").Value = "'5.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.0744"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "'3.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "'16.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.003.95"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "'10.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").Value = "'2.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Value = "'52.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
